$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 123
$ws1.Range("F3").Value = 726
$ws1.Range("F6").Value = 2939
$ws1.Range("F7").Value = 1691
$ws1.Range("F8").Value = 1947
$ws1.Range("F10").Value = 294
$ws1.Range("F11").Value = 795
$ws1.Range("F12").Value = 938
$ws1.Range("F14").Value = 403
$ws1.Range("F19").Value = 7049
$ws1.Range("F21").Value = 1750
$ws1.Range("F22").Value = 184
$ws1.Range("F25").Value = 360
$ws1.Range("F26").Value = 288
$ws1.Range("F28").Value = 1114
$ws1.Range("F29").Value = 936
$ws1.Range("F37").Value = 150
$ws1.Range("F39").Value = 31
$ws1.Range("F40").Value = 147
$ws1.Range("F41").Value = 258

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 123
$ws4.Range("F3").Value = 726
$ws4.Range("F9").Value = 2939
$ws4.Range("F10").Value = 1691
$ws4.Range("F11").Value = 1947
$ws4.Range("F13").Value = 294
$ws4.Range("F14").Value = 795
$ws4.Range("F16").Value = 938
$ws4.Range("F18").Value = 403
$ws4.Range("F22").Value = 7049
$ws4.Range("F24").Value = 1750
$ws4.Range("F26").Value = 184
$ws4.Range("F29").Value = 360
$ws4.Range("F30").Value = 288
$ws4.Range("F32").Value = 1114
$ws4.Range("F33").Value = 936
$ws4.Range("F40").Value = 150
$ws4.Range("F42").Value = 31
$ws4.Range("F43").Value = 147
$ws4.Range("F44").Value = 258
